# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# handback (de-de / zh-cn target files) has completed:
#   - Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#   - Adds "Latest Target File" (F) and "Latest Handback File" (G) hyperlink
#     cells for each data row on the zh-cn and de-de sheets
#   - Updates "Latest Handback DateTime" (H) from the placeholder
#     "0001-01-01 00:00:00" to the real handback timestamps

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the status wording wherever it is used. The text
#    "Ready for handoff" is a shared string also reused (accidentally, by
#    the original report generator) in the "File Extension" column of the
#    zh-cn/de-de sheets, so every cell holding that string must move to the
#    new wording to keep the shared string table / rendered text in sync.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File (F), Latest Handback File (G)
#    and Latest Handback DateTime (H) for both data rows.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/b1e0866879253113d047369139cdd88c607afe36/e2e/4bce5b39-7f43-41ee-8996-7a65256baf22.md",
    "",
    "",
    "4bce5b39-7f43-41ee-8996-7a65256baf22.md"
) | Out-Null

$zhcn.Hyperlinks.Add(
    $zhcn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fe5390b4fef80b3b10d55d15c1537d12138d1932/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/high/4bce5b39-7f43-41ee-8996-7a65256baf22.d1c04fce3085ea57988a1de56ca21a09d325860a.zh-cn.xlf",
    "",
    "",
    "4bce5b39-7f43-41ee-8996-7a65256baf22.d1c04fce3085ea57988a1de56ca21a09d325860a.zh-cn.xlf"
) | Out-Null

$zhcn.Range("H2").Value = "2016-03-18 17:46:01"

$zhcn.Hyperlinks.Add(
    $zhcn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/b1e0866879253113d047369139cdd88c607afe36/e2e/6b285839-d8a4-4d61-af47-499ec345b4b3.md",
    "",
    "",
    "6b285839-d8a4-4d61-af47-499ec345b4b3.md"
) | Out-Null

$zhcn.Hyperlinks.Add(
    $zhcn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fe5390b4fef80b3b10d55d15c1537d12138d1932/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/high/6b285839-d8a4-4d61-af47-499ec345b4b3.ede75e2621103db7a09c4086c642abae97249753.zh-cn.xlf",
    "",
    "",
    "6b285839-d8a4-4d61-af47-499ec345b4b3.ede75e2621103db7a09c4086c642abae97249753.zh-cn.xlf"
) | Out-Null

$zhcn.Range("H3").Value = "2016-03-18 17:46:01"

# ---------------------------------------------------------------------------
# 3. de-de sheet: fill in Latest Target File (F), Latest Handback File (G)
#    and Latest Handback DateTime (H) for both data rows.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/b1e0866879253113d047369139cdd88c607afe36/e2e/4bce5b39-7f43-41ee-8996-7a65256baf22.md",
    "",
    "",
    "4bce5b39-7f43-41ee-8996-7a65256baf22.md"
) | Out-Null

$dede.Hyperlinks.Add(
    $dede.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/76086da8ab2fbe98f4badb3ee4b34962fdf74343/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/high/4bce5b39-7f43-41ee-8996-7a65256baf22.d1c04fce3085ea57988a1de56ca21a09d325860a.de-de.xlf",
    "",
    "",
    "4bce5b39-7f43-41ee-8996-7a65256baf22.d1c04fce3085ea57988a1de56ca21a09d325860a.de-de.xlf"
) | Out-Null

$dede.Range("H2").Value = "2016-03-18 17:46:22"

$dede.Hyperlinks.Add(
    $dede.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/b1e0866879253113d047369139cdd88c607afe36/e2e/6b285839-d8a4-4d61-af47-499ec345b4b3.md",
    "",
    "",
    "6b285839-d8a4-4d61-af47-499ec345b4b3.md"
) | Out-Null

$dede.Hyperlinks.Add(
    $dede.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/76086da8ab2fbe98f4badb3ee4b34962fdf74343/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/high/6b285839-d8a4-4d61-af47-499ec345b4b3.ede75e2621103db7a09c4086c642abae97249753.de-de.xlf",
    "",
    "",
    "6b285839-d8a4-4d61-af47-499ec345b4b3.ede75e2621103db7a09c4086c642abae97249753.de-de.xlf"
) | Out-Null

$dede.Range("H3").Value = "2016-03-18 17:46:22"

Write-Host "Handback report generated."
